$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 377 (404Text) - fill in pending translations
$ws.Range("C377").Value = 'ऐसा लगता है कि यह पेज मौजूद नहीं है'
$ws.Range("E377").Value = 'એવું લાગે છે કે આ પેજ અસ્તિત્વમાં નથી'
$ws.Range("G377").Value = 'মনে হচ্ছে এই পৃষ্ঠাটির অস্তিত্ব নেই'
$ws.Range("I377").Value = 'ఈ పేజీ అందుబాటులో లేదు'
$ws.Range("J377").Value = 'असे दिसते की हे पृष्ठ अस्तित्वात नाही'
$ws.Range("M377").Value = 'ଏହି ପେଜଟି ନଥିଲା ଭଳି ଲାଗୁଛି'
$ws.Range("N377").Value = 'ಈ ಪುಟವು ಅಸ್ತಿತ್ವದಲ್ಲಿಲ್ಲ'

# Row 388 (languageChangeNotification) - update English copy placeholders
$ws.Range("B388").Value = 'You have changed the contribution language from <old language> to <new language>, we will be redirecting you to homepage to start participating again.'

# Row 393 (validationWarningText) - fill in pending translations
$ws.Range("C393").Value = 'हमें लगता है कि आपके द्वारा एन्टर किया गया टेक्स्ट मूल टेक से मेल नहीं खा रहा, क्या आप अपने संपादन के बारे में निश्चित हैं?'
$ws.Range("E393").Value = 'અમને લાગે છે કે તમે દાખલ કરેલ ટેક્સ્ટ મૂળ ટેક્સ્ટ સાથે મેળ ખાતી નથી, શું તમે તમારા સંપાદન વિશે ચોક્કસ છો?'
$ws.Range("G393").Value = 'আমরা মনে করি আপনার লেখা টেক্সট টি মূল টেক্সট এর  সাথে মিলছে না, আপনি কি আপনার সম্পাদনা সম্পর্কে নিশ্চিত ?'
$ws.Range("I393").Value = 'మీరు నమోదు చేసిన వచనం అసలు వచనంతో సరిపోలడం లేదని మేము భావిస్తున్నాము, మీరు మీ సవరణ గురించి ఖచ్చితంగా ఉన్నారా?'
$ws.Range("J393").Value = 'तुम्ही घातलेला मजकूर मूळ मजकुराशी जुळत नाही असे आम्हाला वाटते, तुम्हाला तुमच्या संपादनाबद्दल खात्री आहे का?'
$ws.Range("M393").Value = 'ଆପଣ ଲେଖିଥିବା ଟେକ୍ସଟ ପ୍ରକୃତ ଟେକ୍ସଟ ସହିତ ମେଳ ନ ଖାଇଲା ଭଳିଆଦିଶୁଛି, ଆପଣ ନିଜ ଏଡିଟ ସହିତ ସନ୍ତୁଷ୍ଟ ଅଛନ୍ତି କି ?'
$ws.Range("N393").Value = 'ನೀವು ನಮೂದಿಸಿದ ಪಠ್ಯವು ಮೂಲ ಪಠ್ಯಕ್ಕೆ ಹೊಂದಿಕೆಯಾಗುವುದಿಲ್ಲ ಎಂದು ನಾವು ಭಾವಿಸುತ್ತೇವೆ, ನಿಮ್ಮ ಸಂಪಾದನೆಯ ಬಗ್ಗೆ ನೀವು ಖಚಿತರಾಗಿದ್ದೀರಿಯೇ?'

# New rows 394-406: newly added localization keys
# Row 394: testMic
$ws.Range("A394").Value = 'testMic'
$ws.Range("B394").Value = 'Test mic'
$ws.Range("C394").Value = 'Test Mic'
$ws.Range("E394").Value = 'Test Mic'
$ws.Range("F394").Value = 'Test Mic'
$ws.Range("G394").Value = 'Test Mic'
$ws.Range("H394").Value = 'Test Mic'
$ws.Range("I394").Value = 'Test Mic'
$ws.Range("J394").Value = 'Test Mic'
$ws.Range("K394").Value = 'Test Mic'
$ws.Range("L394").Value = 'Test Mic'
$ws.Range("M394").Value = 'Test Mic'
$ws.Range("N394").Value = 'Test Mic'

# Row 395: playingBackAudio
$ws.Range("A395").Value = 'playingBackAudio'
$ws.Range("B395").Value = 'Playingback Audio'
$ws.Range("C395").Value = 'Playingback Audio'
$ws.Range("E395").Value = 'Playingback Audio'
$ws.Range("F395").Value = 'Playingback Audio'
$ws.Range("G395").Value = 'Playingback Audio'
$ws.Range("H395").Value = 'Playingback Audio'
$ws.Range("I395").Value = 'Playingback Audio'
$ws.Range("J395").Value = 'Playingback Audio'
$ws.Range("K395").Value = 'Playingback Audio'
$ws.Range("L395").Value = 'Playingback Audio'
$ws.Range("M395").Value = 'Playingback Audio'
$ws.Range("N395").Value = 'Playingback Audio'

# Row 396: speakClearly
$ws.Range("A396").Value = 'speakClearly'
$ws.Range("B396").Value = 'Please speak clearly'
$ws.Range("C396").Value = 'Please speak clearly'
$ws.Range("E396").Value = 'Please speak clearly'
$ws.Range("F396").Value = 'Please speak clearly'
$ws.Range("G396").Value = 'Please speak clearly'
$ws.Range("H396").Value = 'Please speak clearly'
$ws.Range("I396").Value = 'Please speak clearly'
$ws.Range("J396").Value = 'Please speak clearly'
$ws.Range("K396").Value = 'Please speak clearly'
$ws.Range("L396").Value = 'Please speak clearly'
$ws.Range("M396").Value = 'Please speak clearly'
$ws.Range("N396").Value = 'Please speak clearly'

# Row 397: recordingCountValidationMsg
$ws.Range("A397").Value = 'recordingCountValidationMsg'
$ws.Range("B397").Value = 'Recording for {{remainingSec}} seconds'
$ws.Range("C397").Value = 'Recording for {{remainingSec}} seconds'
$ws.Range("E397").Value = 'Recording for {{remainingSec}} seconds'
$ws.Range("F397").Value = 'Recording for {{remainingSec}} seconds'
$ws.Range("G397").Value = 'Recording for {{remainingSec}} seconds'
$ws.Range("H397").Value = 'Recording for {{remainingSec}} seconds'
$ws.Range("I397").Value = 'Recording for {{remainingSec}} seconds'
$ws.Range("J397").Value = 'Recording for {{remainingSec}} seconds'
$ws.Range("K397").Value = 'Recording for {{remainingSec}} seconds'
$ws.Range("L397").Value = 'Recording for {{remainingSec}} seconds'
$ws.Range("M397").Value = 'Recording for {{remainingSec}} seconds'
$ws.Range("N397").Value = 'Recording for {{remainingSec}} seconds'

# Row 398: backgroundNoise
$ws.Range("A398").Value = 'backgroundNoise'
$ws.Range("B398").Value = 'Background Noise Detected'
$ws.Range("C398").Value = 'Background Noise Detected'
$ws.Range("E398").Value = 'Background Noise Detected'
$ws.Range("F398").Value = 'Background Noise Detected'
$ws.Range("G398").Value = 'Background Noise Detected'
$ws.Range("H398").Value = 'Background Noise Detected'
$ws.Range("I398").Value = 'Background Noise Detected'
$ws.Range("J398").Value = 'Background Noise Detected'
$ws.Range("K398").Value = 'Background Noise Detected'
$ws.Range("L398").Value = 'Background Noise Detected'
$ws.Range("M398").Value = 'Background Noise Detected'
$ws.Range("N398").Value = 'Background Noise Detected'

# Row 399: lowBackgroundNoise
$ws.Range("A399").Value = 'lowBackgroundNoise'
$ws.Range("B399").Value = 'Low/No Background Noise'
$ws.Range("C399").Value = 'Low/No Background Noise'
$ws.Range("E399").Value = 'Low/No Background Noise'
$ws.Range("F399").Value = 'Low/No Background Noise'
$ws.Range("G399").Value = 'Low/No Background Noise'
$ws.Range("H399").Value = 'Low/No Background Noise'
$ws.Range("I399").Value = 'Low/No Background Noise'
$ws.Range("J399").Value = 'Low/No Background Noise'
$ws.Range("K399").Value = 'Low/No Background Noise'
$ws.Range("L399").Value = 'Low/No Background Noise'
$ws.Range("M399").Value = 'Low/No Background Noise'
$ws.Range("N399").Value = 'Low/No Background Noise'

# Row 400: quickTips
$ws.Range("A400").Value = 'quickTips'
$ws.Range("B400").Value = 'Quick Tips'
$ws.Range("C400").Value = 'Quick Tips'
$ws.Range("E400").Value = 'Quick Tips'
$ws.Range("F400").Value = 'Quick Tips'
$ws.Range("G400").Value = 'Quick Tips'
$ws.Range("H400").Value = 'Quick Tips'
$ws.Range("I400").Value = 'Quick Tips'
$ws.Range("J400").Value = 'Quick Tips'
$ws.Range("K400").Value = 'Quick Tips'
$ws.Range("L400").Value = 'Quick Tips'
$ws.Range("M400").Value = 'Quick Tips'
$ws.Range("N400").Value = 'Quick Tips'

# Row 401: tipOne
$ws.Range("A401").Value = 'tipOne'
$ws.Range("B401").Value = 'Please test your <b>Microphone</b>'
$ws.Range("C401").Value = 'Please test your <b>Microphone</b>'
$ws.Range("E401").Value = 'Please test your <b>Microphone</b>'
$ws.Range("F401").Value = 'Please test your <b>Microphone</b>'
$ws.Range("G401").Value = 'Please test your <b>Microphone</b>'
$ws.Range("H401").Value = 'Please test your <b>Microphone</b>'
$ws.Range("I401").Value = 'Please test your <b>Microphone</b>'
$ws.Range("J401").Value = 'Please test your <b>Microphone</b>'
$ws.Range("K401").Value = 'Please test your <b>Microphone</b>'
$ws.Range("L401").Value = 'Please test your <b>Microphone</b>'
$ws.Range("M401").Value = 'Please test your <b>Microphone</b>'
$ws.Range("N401").Value = 'Please test your <b>Microphone</b>'

# Row 402: tipTwo
$ws.Range("A402").Value = 'tipTwo'
$ws.Range("B402").Value = 'Please test your <b>Speakers</b>'
$ws.Range("C402").Value = 'Please test your <b>Speakers</b>'
$ws.Range("E402").Value = 'Please test your <b>Speakers</b>'
$ws.Range("F402").Value = 'Please test your <b>Speakers</b>'
$ws.Range("G402").Value = 'Please test your <b>Speakers</b>'
$ws.Range("H402").Value = 'Please test your <b>Speakers</b>'
$ws.Range("I402").Value = 'Please test your <b>Speakers</b>'
$ws.Range("J402").Value = 'Please test your <b>Speakers</b>'
$ws.Range("K402").Value = 'Please test your <b>Speakers</b>'
$ws.Range("L402").Value = 'Please test your <b>Speakers</b>'
$ws.Range("M402").Value = 'Please test your <b>Speakers</b>'
$ws.Range("N402").Value = 'Please test your <b>Speakers</b>'

# Row 403: tipThree
$ws.Range("A403").Value = 'tipThree'
$ws.Range("B403").Value = 'Ensure there is <b>no background noise</b>'
$ws.Range("C403").Value = 'Ensure there is <b>no background noise</b>'
$ws.Range("E403").Value = 'Ensure there is <b>no background noise</b>'
$ws.Range("F403").Value = 'Ensure there is <b>no background noise</b>'
$ws.Range("G403").Value = 'Ensure there is <b>no background noise</b>'
$ws.Range("H403").Value = 'Ensure there is <b>no background noise</b>'
$ws.Range("I403").Value = 'Ensure there is <b>no background noise</b>'
$ws.Range("J403").Value = 'Ensure there is <b>no background noise</b>'
$ws.Range("K403").Value = 'Ensure there is <b>no background noise</b>'
$ws.Range("L403").Value = 'Ensure there is <b>no background noise</b>'
$ws.Range("M403").Value = 'Ensure there is <b>no background noise</b>'
$ws.Range("N403").Value = 'Ensure there is <b>no background noise</b>'

# Row 404: tipFour
$ws.Range("A404").Value = 'tipFour'
$ws.Range("B404").Value = 'Read once <b>before recording it</b>'
$ws.Range("C404").Value = 'Read once <b>before recording it</b>'
$ws.Range("E404").Value = 'Read once <b>before recording it</b>'
$ws.Range("F404").Value = 'Read once <b>before recording it</b>'
$ws.Range("G404").Value = 'Read once <b>before recording it</b>'
$ws.Range("H404").Value = 'Read once <b>before recording it</b>'
$ws.Range("I404").Value = 'Read once <b>before recording it</b>'
$ws.Range("J404").Value = 'Read once <b>before recording it</b>'
$ws.Range("K404").Value = 'Read once <b>before recording it</b>'
$ws.Range("L404").Value = 'Read once <b>before recording it</b>'
$ws.Range("M404").Value = 'Read once <b>before recording it</b>'
$ws.Range("N404").Value = 'Read once <b>before recording it</b>'

# Row 405: tipFive
$ws.Range("A405").Value = 'tipFive'
$ws.Range("B405").Value = 'Get started by clicking on <b>Record</b> button'
$ws.Range("C405").Value = 'Get started by clicking on <b>Record</b> button'
$ws.Range("E405").Value = 'Get started by clicking on <b>Record</b> button'
$ws.Range("F405").Value = 'Get started by clicking on <b>Record</b> button'
$ws.Range("G405").Value = 'Get started by clicking on <b>Record</b> button'
$ws.Range("H405").Value = 'Get started by clicking on <b>Record</b> button'
$ws.Range("I405").Value = 'Get started by clicking on <b>Record</b> button'
$ws.Range("J405").Value = 'Get started by clicking on <b>Record</b> button'
$ws.Range("K405").Value = 'Get started by clicking on <b>Record</b> button'
$ws.Range("L405").Value = 'Get started by clicking on <b>Record</b> button'
$ws.Range("M405").Value = 'Get started by clicking on <b>Record</b> button'
$ws.Range("N405").Value = 'Get started by clicking on <b>Record</b> button'

# Row 406: warningAudioPermissionMsg
$ws.Range("A406").Value = 'warningAudioPermissionMsg'
$ws.Range("B406").Value = 'Sorry !!! We could not get access to your audio input device. Make sure you have given microphone access permission'
$ws.Range("C406").Value = 'Sorry !!! We could not get access to your audio input device. Make sure you have given microphone access permission'
$ws.Range("E406").Value = 'Sorry !!! We could not get access to your audio input device. Make sure you have given microphone access permission'
$ws.Range("F406").Value = 'Sorry !!! We could not get access to your audio input device. Make sure you have given microphone access permission'
$ws.Range("G406").Value = 'Sorry !!! We could not get access to your audio input device. Make sure you have given microphone access permission'
$ws.Range("H406").Value = 'Sorry !!! We could not get access to your audio input device. Make sure you have given microphone access permission'
$ws.Range("I406").Value = 'Sorry !!! We could not get access to your audio input device. Make sure you have given microphone access permission'
$ws.Range("J406").Value = 'Sorry !!! We could not get access to your audio input device. Make sure you have given microphone access permission'
$ws.Range("K406").Value = 'Sorry !!! We could not get access to your audio input device. Make sure you have given microphone access permission'
$ws.Range("L406").Value = 'Sorry !!! We could not get access to your audio input device. Make sure you have given microphone access permission'
$ws.Range("M406").Value = 'Sorry !!! We could not get access to your audio input device. Make sure you have given microphone access permission'
$ws.Range("N406").Value = 'Sorry !!! We could not get access to your audio input device. Make sure you have given microphone access permission'

